$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2,3) {
    $ws.Range("K$row").Value = 11
    $ws.Range("M$row").Value = 10.3
    $ws.Range("N$row").Value = 0.04865375531412377
    $ws.Range("O$row").Value = 0.9363636363636364
    $ws.Range("P$row").Value = 10.3
    $ws.Range("Q$row").Value = 0.04865375531412377
    $ws.Range("R$row").Value = 0.9363636363636364
    $ws.Range("U$row").Value = 0.595
    $ws.Range("V$row").Value = 0.002810581010864431
    $ws.Range("W$row").Value = 0.196078431372549
    $ws.Range("X$row").Value = 0.04998170370534875
    $ws.Range("Y$row").Value = 0.1460967276672003
    $ws.Range("AA$row").Value = -0.003204822924541017
    $ws.Range("AB$row").Value = 0.04998170370534875
    $ws.Range("AC$row").Value = -0.05318652662988977
    $ws.Range("AG$row").Value = -0.595
    $ws.Range("AJ$row").Value = -0.00281850264086592
    $ws.Range("AK$row").Value = -0.01175773144946151
    $ws.Range("AM$row").Value = -0.023
    $ws.Range("AQ$row").Value = 8.391304347826088
}
